$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-template_type")

# Add three new rows (122-124) following the same pattern as the
# "reg-ack-template-part3" block (rows 119-121), but for a new
# "reg-ack-template-part4" code.

$ws.Range("A122").Value = "reg-ack-template-part4"
$ws.Range("B122").Value = "Registration Acknowledgement Template - Part 4"
$ws.Range("C122").Value = "eng"
$ws.Range("D122").Value = $true
$ws.Range("E122").Value = "superadmin"
$ws.Range("F122").Value = "now()"

$ws.Range("A123").Value = "reg-ack-template-part4"
$ws.Range("B123").Value = "نموذج شكر التسجيل"
$ws.Range("C123").Value = "ara"
$ws.Range("D123").Value = $true
$ws.Range("E123").Value = "superadmin"
$ws.Range("F123").Value = "now()"

$ws.Range("A124").Value = "reg-ack-template-part4"
$ws.Range("B124").Value = "accusé de réception"
$ws.Range("C124").Value = "fra"
$ws.Range("D124").Value = $true
$ws.Range("E124").Value = "superadmin"
$ws.Range("F124").Value = "now()"

# Update the selection to mirror Excel's behaviour of selecting the
# remainder of the sheet below the newly entered data.
$ws.Range("A125:XFD1048576").Select()
